# Update cryptos list values (price and 1h volume change) to reflect the
# latest scrape, and fix the ordering of TrustWalletToken/TheSandbox rows
# (rows 41-42 swapped name/link/price/volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'23.415.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '

# Row 3
$ws.Range("D3").Value = "'1.630.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.48%  '

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("E5").Value = '  -0.08%  '

# Row 6
$ws.Range("D6").Value = "'304.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '

# Row 7
$ws.Range("D7").Value = "'0.3789"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.66%  '

# Row 8
$ws.Range("D8").Value = "'0.3650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.28%  '

# Row 9
$ws.Range("D9").Value = "'51.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.25%  '

# Row 10
$ws.Range("D10").Value = "'0.08236"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.07%  '

# Row 11
$ws.Range("D11").Value = "'1.235"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.07%  '

# Row 12
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '

# Row 13
$ws.Range("D13").Value = "'22.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.25%  '

# Row 14
$ws.Range("D14").Value = "'6.548"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.16%  '

# Row 15
$ws.Range("D15").Value = "'0.00001249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.89%  '

# Row 16
$ws.Range("D16").Value = "'7.331"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.35%  '

# Row 17
$ws.Range("D17").Value = "'1.629.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.61%  '

# Row 18
$ws.Range("D18").Value = "'94.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '

# Row 19
$ws.Range("D19").Value = "'0.06978"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.69%  '

# Row 20
$ws.Range("D20").Value = "'17.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.36%  '

# Row 21
$ws.Range("D21").Value = "'6.528"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.23%  '

# Row 22
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("D23").Value = "'12.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.42%  '

# Row 24
$ws.Range("D24").Value = "'23.407.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.16%  '

# Row 25
$ws.Range("E25").Value = '  +0.72%  '

# Row 26
$ws.Range("E26").Value = '  +1.15%  '

# Row 27
$ws.Range("D27").Value = "'21.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.70%  '

# Row 28
$ws.Range("D28").Value = "'150.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.12%  '

# Row 29
$ws.Range("D29").Value = "'5.295"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.14%  '

# Row 30
$ws.Range("D30").Value = "'133.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.00%  '

# Row 31
$ws.Range("D31").Value = "'1.806.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.96%  '

# Row 32
$ws.Range("E32").Value = '  -2.12%  '

# Row 33
$ws.Range("D33").Value = "'6.833"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.78%  '

# Row 34
$ws.Range("D34").Value = "'1.030"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.76%  '

# Row 35
$ws.Range("D35").Value = "'10.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.87%  '

# Row 36
$ws.Range("D36").Value = "'0.02786"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.13%  '

# Row 37
$ws.Range("D37").Value = "'0.2523"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.00%  '

# Row 38
$ws.Range("D38").Value = "'0.08785"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.63%  '

# Row 39
$ws.Range("D39").Value = "'0.07113"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.49%  '

# Row 40
$ws.Range("D40").Value = "'6.022"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.7030"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.86%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'1.348"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.94%  '

# Row 43
$ws.Range("D43").Value = "'16.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.46%  '

# Row 44
$ws.Range("D44").Value = "'12.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.00%  '

# Row 45
$ws.Range("D45").Value = "'0.6547"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.26%  '

# Row 47
$ws.Range("D47").Value = "'2.300"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.50%  '

# Row 48
$ws.Range("D48").Value = "'3.977"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.09%  '

# Row 49
$ws.Range("D49").Value = "'0.08019"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.76%  '

# Row 50
$ws.Range("D50").Value = "'1.199"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.42%  '

# Row 51
$ws.Range("D51").Value = "'126.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.13%  '

